# FN-1620: update the utilisation-report fixture so the "fees paid" currency
# columns exist (Fees paid to UKEF currency / Payment currency / Payment
# exchange rate) alongside the existing facility rows, matching the refreshed
# cypress fixture used by the csv/xlsx validators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) -------------------------------------------------
# F1 used to hold "Facility utilisation" via a freshly-appended shared
# string; it now reuses the existing one. G1/H1 stay as they are. Three new
# header cells are appended for the fees-paid-currency columns.
$ws.Range("F1").Value = "Facility utilisation"
$ws.Range("I1").Value = "Fees paid to UKEF currency"
$ws.Range("J1").Value = "Payment currency"
$ws.Range("K1").Value = "Payment exchange rate"

# Copy the existing F1:H1 header formatting onto the new I1:K1 headers so
# they pick up the same green-fill / centered / wrap-text style (style 2).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(1).RowHeight = 121.8

# ---------------------------------------------------------------------------
# 2. Data rows 2-4 --------------------------------------------------------
# Row 2 (Exporter 1 GEF) keeps its existing values; just add the new
# fees-paid-currency columns (matching base currency / GBP, no rate).
$ws.Range("I2").Value = "GBP"
$ws.Range("J2").Value = "GBP"
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("I2:J2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 3 used to be "Exporter 2 GEF" / UKEF 20001372 / EUR / 300000 / ...
# It becomes a second "20001371" facility row for Exporter 2, GBP-based.
$ws.Range("A3").Value = "Exporter 2 GEF"
$ws.Range("B3").Value = 20001371
$ws.Range("C3").Value = "Exporter 2"
$ws.Range("D3").Value = "GBP"
$ws.Range("E3").Value = 600000
$ws.Range("F3").Value = 100000
$ws.Range("G3").Value = 150
$ws.Range("H3").Value = 243
$ws.Range("I3").Value = "GBP"
$ws.Range("J3").Value = "GBP"
# A3/C3 in the source workbook pick up mismatched styles (A3 keeps the usual
# "style 3" text style, but C3 also uses "style 3" rather than "style 5" -
# a leftover from how the row above was edited in Excel).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D2").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 4 used to be "Potato GEF"/"Potato Exporter"/JPY/... and becomes the
# "Potato Gef"/"Potato exporter" GBP row with a EUR fees-paid currency and an
# explicit exchange rate.
$ws.Range("A4").Value = "Potato Gef"
$ws.Range("B4").Value = 20001371
$ws.Range("C4").Value = "Potato exporter"
$ws.Range("D4").Value = "GBP"
$ws.Range("E4").Value = 600000
$ws.Range("F4").Value = 100000
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = "EUR"
$ws.Range("J4").Value = "GBP"
$ws.Range("K4").Value = 1.17

# Row 5 (previously "Fish Exporter") is removed entirely - shift everything
# below up so the old row 5 data disappears and what used to be sparsely
# filled blank rows reappear starting at row 6.
$ws.Rows.Item(5).Delete(-4162) | Out-Null

for ($r = 2; $r -le 4; $r++) {
    $ws.Rows.Item($r).RowHeight = 17.4
}

# ---------------------------------------------------------------------------
# 3. Pre-formatted, empty rows below the data (6-16, then 19-21) --------
# These mirror the column-formatting of the data rows above (copy/paste
# formats only, no values) exactly as Excel leaves behind when a user
# selects/fills a block of rows below their data without typing anything in.
$ws.Range("A4:K4").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:K11").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:K16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A20:H20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B4:H4").Copy() | Out-Null
$ws.Range("B19:H19").PasteSpecial(-4122) | Out-Null
$ws.Range("B21:H21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 14 only has data-style formatting out to column J (I/J only, no K).
$ws.Range("J4").Copy() | Out-Null
$ws.Range("J14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 7 only keeps column A + I:K formatted.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I4:K4").Copy() | Out-Null
$ws.Range("I7:K7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 12 only keeps column K formatted.
$ws.Range("K4").Copy() | Out-Null
$ws.Range("K12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 8 & 13: a lone bold cell in column A (left over heading/note row),
# plus a single "style 5" cell elsewhere.
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").Interior.ColorIndex = 0
$ws.Range("I1").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "Fees paid to UKEF currency"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("K4").Copy() | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null
$ws.Range("J13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($r = 6; $r -le 21; $r++) {
    $ws.Rows.Item($r).RowHeight = 17.4
}

# ---------------------------------------------------------------------------
# 4. Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.39
$ws.Columns.Item(2).ColumnWidth = 16.39
$ws.Columns.Item(3).ColumnWidth = 17.83
$ws.Columns.Item(4).ColumnWidth = 10.17
$ws.Columns.Item(5).ColumnWidth = 15.5
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 14.05
$ws.Columns.Item(8).ColumnWidth = 11.94
$ws.Columns.Item(9).ColumnWidth = 13.28
$ws.Columns.Item(10).ColumnWidth = 12.72
$ws.Columns.Item(11).ColumnWidth = 19.28
$ws.Columns.Item(12).ColumnWidth = 16.39

# ---------------------------------------------------------------------------
# 5. Sheet view / window state ---------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12576

$ws.Range("H2").Select() | Out-Null
$win.ScrollColumn = 2

# ---------------------------------------------------------------------------
# 6. Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
